# Gantt project planner.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Project Presentation (Phase 3 Implementation)" ---
# C6 gains a value (Actual start date), D6 end date pushed out, E6 gets total minutes
$ws.Range("C6").Value = 45763
$ws.Range("D6").Value = 45783
$ws.Range("E6").Value = 1310

# --- Row 7: "Project Presentation (Phase 3 Testing)" ---
# B7 predicted start date moves out, D7 end date cleared
$ws.Range("B7").Value = 45783
$ws.Range("D7").ClearContents()

# --- Row 8: "Project Presentation (Phase 3 Maintenance)" ---
# B8 predicted start date cleared, D8 end date cleared
$ws.Range("B8").ClearContents()
$ws.Range("D8").ClearContents()

# --- Row 53 (table "Table6"): add implementation-phase contribution note for Mahlet Bekele ---
$ws.Range("D53").Value = "server (in the start) , player , Penals, testing"

# --- View settings ---
$ws.Application.ActiveWindow.Zoom = 53
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H13").Select()
